# Regenerate save_data column G ("K") with new computed strikeout values,
# replacing the old "Strike#" based values, per commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 6
    3  = 4
    4  = 4
    5  = 4
    6  = 4
    7  = 5
    8  = 7
    9  = 2
    10 = 7
    11 = 9
    12 = 7
    13 = 2
    14 = 8
    15 = 10
    16 = 4
    17 = 4
    18 = 8
    19 = 7
    20 = 8
    21 = 6
    22 = 0
    23 = 5
    24 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
